$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value looks like a plain number need a leading
# apostrophe (quote-prefix) so Excel stores them as text instead of
# silently converting to a numeric value, matching the original
# inlineStr/text cell type. The style is then reset to "Normal" so the
# quote-prefix formatting flag does not linger on the cell style.

$ws.Range("D2").Value = "34.519.67"
$ws.Range("E2").Value = "  -0.31%  "
$ws.Range("D3").Value = "1.811.73"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "'228.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.42%  "
$ws.Range("E6").Value = "  +3.76%  "
$ws.Range("E7").Value = "  +0.29%  "
$ws.Range("D8").Value = "'34.97"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.62%  "
$ws.Range("E9").Value = "  +1.57%  "
$ws.Range("E10").Value = "  -0.18%  "
$ws.Range("D12").Value = "2.073.41"
$ws.Range("E12").Value = "  +0.55%  "
$ws.Range("D13").Value = "'11.24"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.81%  "
$ws.Range("D14").Value = "1.801.59"
$ws.Range("E14").Value = "  +0.07%  "
$ws.Range("D15").Value = "'0.647"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.14%  "
$ws.Range("D16").Value = "34.517.48"
$ws.Range("E16").Value = "  -0.24%  "
$ws.Range("E17").Value = "  +2.60%  "
$ws.Range("D18").Value = "'69.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.42%  "
$ws.Range("E19").Value = "  -0.86%  "
$ws.Range("D20").Value = "'245.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.84%  "
$ws.Range("D21").Value = "'11.45"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.88%  "
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("E23").Value = "  -0.72%  "
$ws.Range("D24").Value = "'172.71"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.65%  "
$ws.Range("D25").Value = "'2.11"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.74%  "
$ws.Range("D26").Value = "'7.96"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +8.75%  "
$ws.Range("E27").Value = "  +1.24%  "
$ws.Range("E28").Value = "  +2.12%  "
$ws.Range("E30").Value = "  -2.48%  "
$ws.Range("E31").Value = "  +1.19%  "
$ws.Range("E32").Value = "  +0.94%  "
$ws.Range("D33").Value = "'1.24"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("E34").Value = "  -0.27%  "
$ws.Range("B35").Value = "Maker"
$ws.Range("C35").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D35").Value = "1.394.88"
$ws.Range("E35").Value = "  -2.59%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'0.681"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.47%  "
$ws.Range("D37").Value = "'2.46"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.36%  "
$ws.Range("E38").Value = "  -1.20%  "
$ws.Range("E39").Value = "  -0.48%  "
$ws.Range("D40").Value = "'83.84"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.79%  "
$ws.Range("E41").Value = "  +1.45%  "
$ws.Range("E42").Value = "  +2.64%  "
$ws.Range("E43").Value = "  -0.78%  "
$ws.Range("D44").Value = "'13.36"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.47%  "
$ws.Range("E45").Value = "  +4.05%  "
$ws.Range("E46").Value = "  -1.97%  "
$ws.Range("D47").Value = "'5.99"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.02%  "
$ws.Range("D48").Value = "1.973.19"
$ws.Range("E48").Value = "  +0.63%  "
$ws.Range("D49").Value = "'105.18"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.62%  "
$ws.Range("E50").Value = "  +2.61%  "
$ws.Range("E51").Value = "  +0.22%  "
